{"js": "// Add the missing description and input-info text to the \"(5) Xem quy\u1ec1n\"\n// (View permissions) section: append a space run + the description text\n// to the \"- M\u00f4 t\u1ea3:\" paragraph, and append a space run to the\n// \"- Th\u00f4ng tin \u0111\u1ea7u v\u00e0o:\" paragraph \u2014 each as its own new run (matching\n// how the other similar sections in the document already look).\n\n// Build a minimal \"Flat OPC\" package wrapping a <w:p> fragment containing\n// the given run(s) XML, suitable for Paragraph.insertOoxml/Range.insertOoxml.\nfunction flatOpcParagraph(runsXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + runsXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n}\n\nfunction runXml(text, preserveSpace) {\n  const sp = preserveSpace ? ' xml:space=\"preserve\"' : '';\n  return '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t' + sp + '>' + text + '</w:t></w:r>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two target paragraphs: the \"- M\u00f4 t\u1ea3:\" and the\n// \"- Th\u00f4ng tin \u0111\u1ea7u v\u00e0o:\" paragraph that immediately follow\n// \"(5) Xem quy\u1ec1n:\" (they are currently empty / label-only).\nlet descriptionPara = null;\nlet inputInfoPara = null;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = (items[i].text || \"\").trim();\n  if (t === \"(5) Xem quy\u1ec1n:\") {\n    descriptionPara = items[i + 1];\n    inputInfoPara = items[i + 2];\n    break;\n  }\n}\n\nif (!descriptionPara || !inputInfoPara) {\n  throw new Error(\"Could not locate the '(5) Xem quy\u1ec1n' description/input paragraphs.\");\n}\n\n// Paragraph 1 (\"- M\u00f4 t\u1ea3:\") gets a space run followed by the description run.\nconst descRunsXml = runXml(\" \", true) +\n  runXml(\"H\u1ec7 th\u1ed1ng cho ph\u00e9p ng\u01b0\u1eddi qu\u1ea3n tr\u1ecb h\u1ec7 th\u1ed1ng xem th\u00f4ng tin quy\u1ec1n\", false);\ndescriptionPara.insertOoxml(flatOpcParagraph(descRunsXml), \"End\");\n\n// Paragraph 2 (\"- Th\u00f4ng tin \u0111\u1ea7u v\u00e0o:\") gets a trailing space run.\nconst inputRunsXml = runXml(\" \", true);\ninputInfoPara.insertOoxml(flatOpcParagraph(inputRunsXml), \"End\");\n\nawait context.sync();\n", "ps1": "# Add the missing description and input-info text to the \"(5) Xem quyen\"\n# (View permissions) section: append a space run + the description text\n# to the \"- Mo ta:\" paragraph, and append a space run to the\n# \"- Thong tin dau vao:\" paragraph - each as its own new run (matching\n# how the other similar sections in the document already look).\n#\n# (ASCII used in comments only; the actual inserted document text below\n# uses the correct Vietnamese diacritics.)\n\n$d = $word.ActiveDocument\n\n# Locate the \"(5) Xem quy\u1ec1n:\" heading paragraph; the description and\n# input-info paragraphs immediately follow it.\n$headingIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.Trim()\n    if ($t -eq \"(5) Xem quy\u1ec1n:\") {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -eq -1) {\n    throw \"Could not locate the '(5) Xem quy\u1ec1n:' heading paragraph.\"\n}\n\n$descPara = $d.Paragraphs($headingIndex + 1)\n$inputPara = $d.Paragraphs($headingIndex + 2)\n\nfunction New-RunFlatOpc($runsXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\nfunction New-RunXml($text, [bool]$preserveSpace) {\n    $sp = \"\"\n    if ($preserveSpace) { $sp = ' xml:space=\"preserve\"' }\n    return '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t' + $sp + '>' + $text + '</w:t></w:r>'\n}\n\n# Paragraph 1 (\"- M\u00f4 t\u1ea3:\") gets a space run followed by the description run.\n$descRange = $descPara.Range\n$descRange.End = $descRange.End - 1   # position right before the paragraph mark\n$descRange.Collapse(0)                # wdCollapseEnd\n$descRunsXml = (New-RunXml \" \" $true) + (New-RunXml \"H\u1ec7 th\u1ed1ng cho ph\u00e9p ng\u01b0\u1eddi qu\u1ea3n tr\u1ecb h\u1ec7 th\u1ed1ng xem th\u00f4ng tin quy\u1ec1n\" $false)\n$descRange.InsertXML((New-RunFlatOpc $descRunsXml), \"End\")\n\n# Paragraph 2 (\"- Th\u00f4ng tin \u0111\u1ea7u v\u00e0o:\") gets a trailing space run.\n$inputRange = $inputPara.Range\n$inputRange.End = $inputRange.End - 1\n$inputRange.Collapse(0)\n$inputRunsXml = New-RunXml \" \" $true\n$inputRange.InsertXML((New-RunFlatOpc $inputRunsXml), \"End\")\n"}
